$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Promotions" column data (C1 header text, and C2:C4 values)
$ws.Range("C1").ClearContents()
$ws.Range("C2:C4").ClearContents()
